# "Generate Report for Handoff"
#
# Refreshes the localization-status report: the three files that were
# "Ready for handoff" (3ef57ec2, 6d3486c6, acd663b6, e7857b38) get their
# Priority bumped from "low" to "ht", and the Latest Handoff Datetime for
# both target languages moves forward a few seconds to reflect the new
# handoff run. The Overview sheet's "Latest HO Xliff Generate Date" column
# mirrors the de-de handoff timestamp, so it is refreshed to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 (the "Ready for handoff" files)
#   Priority (E): low -> ht
#   Latest Handoff Datetime (H): 2016-08-30 02:33:10 -> 2016-08-30 02:33:26
foreach ($row in 4..7) {
    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-08-30 02:33:26"
}

# de-de sheet: rows 4-7 (the "Ready for handoff" files)
#   Priority (E): low -> ht
#   Latest Handoff Datetime (H): 2016-08-30 02:33:15 -> 2016-08-30 02:33:31
foreach ($row in 4..7) {
    $wsDeDe.Range("E$row").Value = "ht"
    $wsDeDe.Range("H$row").Value = "2016-08-30 02:33:31"
}

# Overview sheet: rows 4-7, Latest HO Xliff Generate Date (G) mirrors the
# de-de handoff timestamp above.
foreach ($row in 4..7) {
    $wsOverview.Range("G$row").Value = "2016-08-30 02:33:31"
}
